$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 39011.168
$ws.Range("J47").Value = 41000
$ws.Range("L47").Value = 41000
$ws.Range("N47").Value = -42944
$ws.Range("H93").Value = 47950.75
$ws.Range("J93").Value = 47950.75
$ws.Range("L93").Value = 47950.75
$ws.Range("N93").Value = -52942.75
$ws.Range("H106").Value = 848.5454999999999
$ws.Range("I106").Value = 791.75
$ws.Range("J106").Value = 1000
$ws.Range("K106").Value = 791.75
$ws.Range("L106").Value = 1000
$ws.Range("M106").Value = -160.75
$ws.Range("N106").Value = -2262
$ws.Range("H112").Value = 1970.0667
$ws.Range("I112").Value = 800.2
$ws.Range("J112").Value = 2204.04
$ws.Range("K112").Value = 2400.6
$ws.Range("L112").Value = 6612.12
$ws.Range("M112").Value = -1292.6
$ws.Range("N112").Value = -8828.119999999999
$ws.Range("H113").Value = 3861
$ws.Range("I113").Value = 3826.25
$ws.Range("K113").Value = 3826.25
$ws.Range("M113").Value = -572.25
$ws.Range("H116").Value = 33336100
$ws.Range("I116").Value = 66668800
$ws.Range("J116").Value = 3400
$ws.Range("K116").Value = 66668800
$ws.Range("L116").Value = 3400
$ws.Range("M116").Value = -66665358
$ws.Range("N116").Value = -10284
$ws.Range("H133").Value = 34992.5
$ws.Range("J133").Value = 34992.5
$ws.Range("L133").Value = 34992.5
$ws.Range("N133").Value = -45112.5
$ws.Range("H138").Value = 2146.74
$ws.Range("I138").Value = 1312
$ws.Range("J138").Value = 2521.768
$ws.Range("K138").Value = 3936
$ws.Range("L138").Value = 7565.304
$ws.Range("M138").Value = 1204
$ws.Range("N138").Value = -17845.304
$ws.Range("H139").Value = 58923.8
$ws.Range("J139").Value = 58923.8
$ws.Range("L139").Value = 58923.8
$ws.Range("N139").Value = -69203.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2513.7144
$ws.Range("I110").Value = 2479.2
$ws.Range("J110").Value = 2600
$ws.Range("K110").Value = 2479.2
$ws.Range("L110").Value = 2600
$ws.Range("M110").Value = -434.1999999999998
$ws.Range("N110").Value = -6690
$ws.Range("H132").Value = 1118686.2
$ws.Range("I132").Value = 1252834.8
$ws.Range("J132").Value = 45499
$ws.Range("K132").Value = 3758504.4
$ws.Range("L132").Value = 136497
$ws.Range("M132").Value = -3755974.4
$ws.Range("N132").Value = -141557

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 59999.5
$ws.Range("J95").Value = 59999.5
$ws.Range("L95").Value = 59999.5
$ws.Range("N95").Value = -65491.5
$ws.Range("H107").Value = 28516.895
$ws.Range("I107").Value = 33681.688
$ws.Range("J107").Value = 971.3333
$ws.Range("K107").Value = 33681.688
$ws.Range("L107").Value = 971.3333
$ws.Range("M107").Value = -31761.688
$ws.Range("N107").Value = -4811.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8666.666999999999
$ws.Range("J4").Value = 8666.666999999999
$ws.Range("L4").Value = 8666.666999999999
$ws.Range("N4").Value = -8890.666999999999
$ws.Range("H32").Value = 3106.4
$ws.Range("I32").Value = 503.33334
$ws.Range("J32").Value = 7011
$ws.Range("K32").Value = 503.33334
$ws.Range("L32").Value = 7011
$ws.Range("M32").Value = -187.33334
$ws.Range("N32").Value = -7643
$ws.Range("H63").Value = 23415.334
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H66").Value = 23415.334
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H68").Value = 31765
$ws.Range("J68").Value = 31765
$ws.Range("L68").Value = 31765
$ws.Range("N68").Value = -33263
$ws.Range("H71").Value = 31765
$ws.Range("J71").Value = 31765
$ws.Range("L71").Value = 95295
$ws.Range("N71").Value = -102783
$ws.Range("H74").Value = 32140
$ws.Range("J74").Value = 32140
$ws.Range("L74").Value = 32140
$ws.Range("N74").Value = -33888
$ws.Range("H77").Value = 32140
$ws.Range("J77").Value = 32140
$ws.Range("L77").Value = 96420
$ws.Range("N77").Value = -105156
$ws.Range("H81").Value = 35666.668
$ws.Range("J81").Value = 35666.668
$ws.Range("L81").Value = 35666.668
$ws.Range("N81").Value = -37662.668
$ws.Range("H84").Value = 35666.668
$ws.Range("J84").Value = 35666.668
$ws.Range("L84").Value = 107000.004
$ws.Range("N84").Value = -116984.004
$ws.Range("H93").Value = 13083.917
$ws.Range("I93").Value = 8172.4287
$ws.Range("J93").Value = 19960
$ws.Range("K93").Value = 8172.4287
$ws.Range("L93").Value = 19960
$ws.Range("M93").Value = -6300.4287
$ws.Range("N93").Value = -23704
$ws.Range("H96").Value = 40800
$ws.Range("J96").Value = 40800
$ws.Range("L96").Value = 40800
$ws.Range("N96").Value = -46292
$ws.Range("H104").Value = 20285
$ws.Range("J104").Value = 20285
$ws.Range("L104").Value = 20285
$ws.Range("N104").Value = -25527
$ws.Range("H110").Value = 30719.2
$ws.Range("J110").Value = 30719.2
$ws.Range("L110").Value = 30719.2
$ws.Range("N110").Value = -38899.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 8000
$ws.Range("J47").Value = 8000
$ws.Range("L47").Value = 8000
$ws.Range("N47").Value = -9136
$ws.Range("H75").Value = 46666.668
$ws.Range("J75").Value = 46666.668
$ws.Range("L75").Value = 46666.668
$ws.Range("N75").Value = -48414.668
$ws.Range("H78").Value = 46666.668
$ws.Range("J78").Value = 46666.668
$ws.Range("L78").Value = 140000.004
$ws.Range("N78").Value = -148736.004
$ws.Range("H97").Value = 58676.668
$ws.Range("I97").Value = 85931.664
$ws.Range("J97").Value = 4166.6665
$ws.Range("K97").Value = 85931.664
$ws.Range("L97").Value = 4166.6665
$ws.Range("M97").Value = -85435.664
$ws.Range("N97").Value = -5158.6665
$ws.Range("H122").Value = 4507.2856
$ws.Range("I122").Value = 4343.048
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 13029.144
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -10579.144
$ws.Range("N122").Value = -19900
$ws.Range("H126").Value = 3098.625
$ws.Range("I126").Value = 2237.7144
$ws.Range("J126").Value = 3768.2222
$ws.Range("K126").Value = 6713.1432
$ws.Range("L126").Value = 11304.6666
$ws.Range("M126").Value = -4243.1432
$ws.Range("N126").Value = -16244.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4599.8
$ws.Range("I7").Value = 4999.75
$ws.Range("K7").Value = 4999.75
$ws.Range("M7").Value = -4887.75
$ws.Range("H40").Value = 5964
$ws.Range("I40").Value = 5868.6924
$ws.Range("J40").Value = 6141
$ws.Range("K40").Value = 5868.6924
$ws.Range("L40").Value = 6141
$ws.Range("M40").Value = -5732.6924
$ws.Range("N40").Value = -6413
$ws.Range("H97").Value = 24938.092
$ws.Range("J97").Value = 24938.092
$ws.Range("L97").Value = 24938.092
$ws.Range("N97").Value = -26920.092
$ws.Range("H122").Value = 5639.1377
$ws.Range("I122").Value = 6072.3687
$ws.Range("J122").Value = 4816
$ws.Range("K122").Value = 18217.1061
$ws.Range("L122").Value = 14448
$ws.Range("M122").Value = -15767.1061
$ws.Range("N122").Value = -19348
$ws.Range("H126").Value = 4599.8
$ws.Range("I126").Value = 4999.75
$ws.Range("K126").Value = 14999.25
$ws.Range("M126").Value = -12529.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 93827.17999999999
$ws.Range("I81").Value = 168449.83
$ws.Range("J81").Value = 4280
$ws.Range("K81").Value = 336899.66
$ws.Range("L81").Value = 8560
$ws.Range("M81").Value = -335838.66
$ws.Range("N81").Value = -10682
$ws.Range("H84").Value = 93827.17999999999
$ws.Range("I84").Value = 168449.83
$ws.Range("J84").Value = 4280
$ws.Range("K84").Value = 1684498.3
$ws.Range("L84").Value = 42800
$ws.Range("M84").Value = -1679194.3
$ws.Range("N84").Value = -53408
$ws.Range("H100").Value = 19309.363
$ws.Range("I100").Value = 66833.336
$ws.Range("J100").Value = 1487.875
$ws.Range("K100").Value = 133666.672
$ws.Range("L100").Value = 2975.75
$ws.Range("M100").Value = -133125.672
$ws.Range("N100").Value = -4057.75
$ws.Range("H126").Value = 8479.546
$ws.Range("I126").Value = 9197.5
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 27592.5
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = -25122.5
$ws.Range("N126").Value = -8840
